$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 5000
$ws.Range("M3").Value = 2000
$ws.Range("M4").Value = 3000

$ws.Range("L5").Select()
